$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$srcRow = 15
$dstRow = 16

# Copy formatting (styles/number formats) from the row above so the new
# row matches the existing rows without introducing new style entries.
$ws.Range("A" + $srcRow + ":N" + $srcRow).Copy()
$ws.Range("A" + $dstRow + ":N" + $dstRow).PasteSpecial(-4122)

$ws.Cells.Item($dstRow, 1).Value = 42622.890451388892
$ws.Cells.Item($dstRow, 2).Value = 12
$ws.Cells.Item($dstRow, 3).Value = 62
$ws.Cells.Item($dstRow, 4).Value = 35
$ws.Cells.Item($dstRow, 5).Value = 62
$ws.Cells.Item($dstRow, 6).Value = 24
$ws.Cells.Item($dstRow, 7).Value = 20402
$ws.Cells.Item($dstRow, 8).Value = 28943
$ws.Cells.Item($dstRow, 9).Value = 3220
$ws.Cells.Item($dstRow, 10).Value = 415
$ws.Cells.Item($dstRow, 11).Value = 239
$ws.Cells.Item($dstRow, 12).Value = 44
$ws.Cells.Item($dstRow, 13).Value = 14
$ws.Cells.Item($dstRow, 14).Value = "Bag"
